# Apply the FlashScore 2024-11-18 weekly-games data refresh.
#
# The diff shows the match that was at row 3 (ARGENTINA - TORNEO BETANO,
# Instituto vs Argentinos Jrs, id W2Rn64T7) was removed, so the fixtures that
# were in rows 4-6 shift up to rows 3-5 (dimension goes from A1:BD6 to A1:BD5).
# On top of that row-shift, a number of odds cells across rows 2-5 were
# refreshed with new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 3 (old "W2Rn64T7" / Instituto vs Argentinos Jrs fixture).
# This shifts old rows 4,5,6 up to become rows 3,4,5 and updates the used-range dimension to A1:BD5.
$ws.Rows("3").Delete()

# Step 2: apply the odds tweaks for each (now-shifted) row.

# Row 2
$ws.Cells.Item(2, 8).Value = 3
$ws.Cells.Item(2, 9).Value = 3.8
$ws.Cells.Item(2, 10).Value = 3
$ws.Cells.Item(2, 13).Value = 1.11
$ws.Cells.Item(2, 14).Value = 6.5
$ws.Cells.Item(2, 27).Value = 21
$ws.Cells.Item(2, 35).Value = 15
$ws.Cells.Item(2, 49).Value = 5.5
$ws.Cells.Item(2, 50).Value = 23

# Row 3
$ws.Cells.Item(3, 7).Value = 2.55
$ws.Cells.Item(3, 9).Value = 3
$ws.Cells.Item(3, 10).Value = 3.5
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 15).Value = 1.57
$ws.Cells.Item(3, 16).Value = 2.25
$ws.Cells.Item(3, 19).Value = 1.62
$ws.Cells.Item(3, 20).Value = 2.2
$ws.Cells.Item(3, 21).Value = 2.25
$ws.Cells.Item(3, 22).Value = 1.57
$ws.Cells.Item(3, 23).Value = 6
$ws.Cells.Item(3, 24).Value = 11
$ws.Cells.Item(3, 26).Value = 26
$ws.Cells.Item(3, 29).Value = 6
$ws.Cells.Item(3, 33).Value = 6.5
$ws.Cells.Item(3, 34).Value = 13
$ws.Cells.Item(3, 35).Value = 12
$ws.Cells.Item(3, 36).Value = 34
$ws.Cells.Item(3, 37).Value = 29
$ws.Cells.Item(3, 38).Value = 41
$ws.Cells.Item(3, 40).Value = 4.33
$ws.Cells.Item(3, 41).Value = 17
$ws.Cells.Item(3, 46).Value = 2.2
$ws.Cells.Item(3, 47).Value = 9.5
$ws.Cells.Item(3, 49).Value = 4.75
$ws.Cells.Item(3, 50).Value = 19
$ws.Cells.Item(3, 51).Value = 34
$ws.Cells.Item(3, 52).Value = 67
$ws.Cells.Item(3, 53).Value = 126
$ws.Cells.Item(3, 54).Value = 351

# Row 4
$ws.Cells.Item(4, 7).Value = 1.4
$ws.Cells.Item(4, 9).Value = 7.5
$ws.Cells.Item(4, 10).Value = 1.91
$ws.Cells.Item(4, 12).Value = 7
$ws.Cells.Item(4, 14).Value = 12
$ws.Cells.Item(4, 21).Value = 2
$ws.Cells.Item(4, 22).Value = 1.73
$ws.Cells.Item(4, 24).Value = 6.5
$ws.Cells.Item(4, 26).Value = 9
$ws.Cells.Item(4, 30).Value = 9
$ws.Cells.Item(4, 31).Value = 21
$ws.Cells.Item(4, 33).Value = 19
$ws.Cells.Item(4, 34).Value = 41
$ws.Cells.Item(4, 39).Value = 351
$ws.Cells.Item(4, 43).Value = 19
$ws.Cells.Item(4, 49).Value = 8.5
$ws.Cells.Item(4, 50).Value = 41
$ws.Cells.Item(4, 52).Value = 151

# Row 5
$ws.Cells.Item(5, 7).Value = 1.95
$ws.Cells.Item(5, 8).Value = 3.1
$ws.Cells.Item(5, 9).Value = 4.33
$ws.Cells.Item(5, 10).Value = 2.63
$ws.Cells.Item(5, 12).Value = 4.75
$ws.Cells.Item(5, 21).Value = 2
$ws.Cells.Item(5, 22).Value = 1.73
$ws.Cells.Item(5, 23).Value = 6
$ws.Cells.Item(5, 24).Value = 8.5
$ws.Cells.Item(5, 25).Value = 9
$ws.Cells.Item(5, 26).Value = 17
$ws.Cells.Item(5, 27).Value = 17
$ws.Cells.Item(5, 31).Value = 17
$ws.Cells.Item(5, 34).Value = 21
$ws.Cells.Item(5, 35).Value = 15
$ws.Cells.Item(5, 37).Value = 41
$ws.Cells.Item(5, 39).Value = 351
$ws.Cells.Item(5, 40).Value = 3.75
$ws.Cells.Item(5, 41).Value = 11
$ws.Cells.Item(5, 47).Value = 9
$ws.Cells.Item(5, 49).Value = 6
$ws.Cells.Item(5, 50).Value = 23
$ws.Cells.Item(5, 51).Value = 34
$ws.Cells.Item(5, 52).Value = 81
$ws.Cells.Item(5, 53).Value = 126
$ws.Cells.Item(5, 54).Value = 301
